$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 59, pushing existing rows 59-152
# down to 61-154 (formatting/styles of the surrounding rows are inherited
# automatically by Insert()).
$ws.Range("A59:A60").EntireRow.Insert()

# New row 59: Pepino ensalada - "Primera" entry for the newest reporting date.
$ws.Range("A59").Value = 9
$ws.Range("B59").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C59").Value = "Metropolitana"
$ws.Range("D59").Value = (Get-Date -Year 2021 -Month 9 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E59").Value = 13
$ws.Range("F59").Value = 100112043
$ws.Range("G59").Value = "Pepino ensalada"
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 133
$ws.Range("K59").Value = 14000
$ws.Range("L59").Value = 15000
$ws.Range("M59").Value = 14496
$ws.Range("N59").Value = "`$/caja 60 unidades"
$ws.Range("O59").Value = "Región de Arica y Parinacota"
$ws.Range("P59").Value = 242
$ws.Range("Q59").Value = 60
$ws.Range("R59").Value = "Hortaliza"

# New row 60: Pepino ensalada - "Segunda" entry for the newest reporting date.
$ws.Range("A60").Value = 9
$ws.Range("B60").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C60").Value = "Metropolitana"
$ws.Range("D60").Value = (Get-Date -Year 2021 -Month 9 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E60").Value = 13
$ws.Range("F60").Value = 100112043
$ws.Range("G60").Value = "Pepino ensalada"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Segunda"
$ws.Range("J60").Value = 79
$ws.Range("K60").Value = 12000
$ws.Range("L60").Value = 13000
$ws.Range("M60").Value = 12494
$ws.Range("N60").Value = "`$/caja 100 unidades"
$ws.Range("O60").Value = "Región de Arica y Parinacota"
$ws.Range("P60").Value = 125
$ws.Range("Q60").Value = 100
$ws.Range("R60").Value = "Hortaliza"
